# Apply the "outputToCloud(resource)" / "text" category addition to the
# hidden '#system' lookup sheet, matching the authoring commit:
#   [base] - [`outputToCloud(resource)`]: support the transferring of
#   output artifact to the cloud.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Make room for a brand new "text" category column. It sorts
#    alphabetically right before "web", so insert a fresh column in
#    front of the existing column Y (which shifts Y..AD to Z..AE).
# ---------------------------------------------------------------------
$ws.Columns("Y").Insert()

$ws.Cells.Item(1, 25).Value2 = "text"
$ws.Cells.Item(2, 25).Value2 = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------
# 2) Insert the new "text" category name into the sorted "target" list
#    (column A), which lives right before "web" (was row 25, now 26).
# ---------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $ws.Cells.Item($r + 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
}
$ws.Cells.Item(25, 1).Value2 = "text"

# ---------------------------------------------------------------------
# 3) Insert the new "outputToCloud(resource)" function into the sorted
#    "base" function list (column E), right before "prependText(...)".
# ---------------------------------------------------------------------
for ($r = 38; $r -ge 22; $r--) {
    $ws.Cells.Item($r + 1, 5).Value2 = $ws.Cells.Item($r, 5).Value2
}
$ws.Cells.Item(22, 5).Value2 = "outputToCloud(resource)"

# ---------------------------------------------------------------------
# 4) Refresh the named ranges so they point at the new boundaries.
#    (Named ranges are not auto-adjusted by the column insert above.)
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"

$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
